# Commit: "Remove wrong holidays date in example 39"
#
# The "holiday" sheet listed 2025-01-01 (serial 45658) as a holiday, which
# was wrong, so that row is removed. Removing it shifts the remaining
# holiday dates up by one row and automatically updates every
# NETWORKDAYS(..., holiday!A$2:A$500) reference on the "task" sheet to
# holiday!A$2:A$499 (and recalculates all dependent values across the
# workbook, e.g. the "xbday" sheet's VLOOKUP results).
#
# The row delete is also what caused the previously-selected/active sheet
# ("holiday") to become the active tab instead of "misc".

$wb = $excel.ActiveWorkbook

$holidayWs = $wb.Worksheets.Item("holiday")

# Remove the erroneous holiday date (row 2, serial 45658 = 2025-01-01).
# This shifts subsequent rows up and automatically adjusts every formula
# that referenced the holiday!A$2:A$500 range (e.g. on the "task" sheet)
# down to holiday!A$2:A$499.
$holidayWs.Rows.Item(2).Delete()

# Recalculate the whole workbook so every dependent value (task sheet
# NETWORKDAYS results, ratios, floor/ceiling helpers, xbday VLOOKUPs, ...)
# is refreshed.
$wb.Application.Calculate()

# Make "holiday" the active sheet/tab, with the same selected cell the
# author ended up with.
$holidayWs.Activate()
$holidayWs.Range("E12").Select()
